$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (orig_acc)
$ws.Range("B2").Value = 0.9020000696182251
$ws.Range("C2").Value = 0.8020000457763672
$ws.Range("D2").Value = 0.7950000166893005
$ws.Range("E2").Value = 0.8530000448226929

# Row 3 (orig_sim)
$ws.Range("C3").Value = 0.8450000286102295
$ws.Range("D3").Value = 0.8330000638961792
$ws.Range("E3").Value = 0.8860000371932983

# Row 4 (orig_acc_robust)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.7994269728660583
$ws.Range("D4").Value = 0.7939913868904114
$ws.Range("E4").Value = 0.8428571224212646

# Row 5 (orig_sim_robust)
$ws.Range("C5").Value = 0.830945611000061
$ws.Range("D5").Value = 0.8283261656761169
$ws.Range("E5").Value = 0.8285714387893677

# Row 6 (orig_acc_adv)
$ws.Range("B6").Value = 0.9018036127090454
$ws.Range("C6").Value = 0.8033794164657593
$ws.Range("D6").Value = 0.795306384563446
$ws.Range("E6").Value = 0.8537634611129761

# Row 7 (orig_sim_adv)
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0.8525345325469971
$ws.Range("D7").Value = 0.8344197869300842
$ws.Range("E7").Value = 0.8903226256370544

# Row 8 (adv_sim)
$ws.Range("C8").Value = 0.6500000357627869
$ws.Range("D8").Value = 0.7660000324249268
$ws.Range("E8").Value = 0.9290000200271606

# Row 9 (adv_hit)
$ws.Range("B9").Value = 0.9980000257492065
$ws.Range("C9").Value = 0.6510000228881836
$ws.Range("D9").Value = 0.7670000195503235
$ws.Range("E9").Value = 0.9300000667572021

# Row 10 (g_align)
$ws.Range("C10").Value = 1.069687604904175
$ws.Range("D10").Value = 1.119789004325867
$ws.Range("E10").Value = 0.9930899143218994

# Row 11 (g_align_robust)
$ws.Range("B11").Value = 0.0006905339541845024
$ws.Range("C11").Value = 1.085886120796204
$ws.Range("D11").Value = 1.204179525375366
$ws.Range("E11").Value = 1.061654448509216

# Row 12 (g_align_adv)
$ws.Range("C12").Value = 1.069687604904175
$ws.Range("D12").Value = 1.119789004325867
$ws.Range("E12").Value = 0.9930899143218994
